# Append a new data row (row 49) to every data sheet in the workbook.
# The new row duplicates the values of the last existing row (row 48)
# except for column A ("time"), which gets the new timestamp recorded
# for this upload.

$wb = $excel.ActiveWorkbook

$newTimestamp = 45835.46310185185

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $srcRow = 48
    $dstRow = 49

    # Column A: new upload timestamp, keep the same date/time number format
    # used by the existing rows.
    $ws.Cells.Item($dstRow, 1).Value = $newTimestamp
    $ws.Cells.Item($dstRow, 1).NumberFormat = $ws.Cells.Item($srcRow, 1).NumberFormat()

    # Columns B through I: carry the previous row's values forward unchanged.
    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($dstRow, $col).Value = $ws.Cells.Item($srcRow, $col).Value()
    }
}
